$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the BAT component row (row 2) and the CONN1 component row (row 6)
# by clearing their cell contents (leaving the row slots empty, not shifting
# subsequent rows up).
$ws.Range("A2:J2").ClearContents()
$ws.Range("A6:J6").ClearContents()
